$wb = $excel.ActiveWorkbook

# --- Sheet "Settings" --------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# New description text in column C for the URL-related settings rows.
$settings.Range("C5").Value = "Source URL of boys names for processing."
$settings.Range("C6").Value = "Source URL of girls names for processing."
$settings.Range("C7").Value = "URL of Unicorn Name generator."

# NamesAmount default value changed 7 -> 10, plus its description.
$settings.Range("B9").Value = 10
$settings.Range("C9").Value = "Determines how many names of each gender should be processed, starting from most popular."

# Move the active selection to B13.
$settings.Activate()
$settings.Range("B13").Select()

# --- Sheet "Constants" ---------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

# MaxRetryNumber value changed 0 -> 1.
$constants.Range("B2").Value = 1

# New rows describing Input/Template/Output paths.
$constants.Range("A12").Value = "InputPath"
$constants.Range("B12").Value = "Data\Input"
$constants.Range("C12").Value = "Path for input files. Could be relative or full path."

$constants.Range("A13").Value = "TemplateFileName"
$constants.Range("B13").Value = "template.xlsx"
$constants.Range("C13").Value = "Filename of template excel spreadsheet."

$constants.Range("A14").Value = "OutputPath"
$constants.Range("B14").Value = "Data\Output"
$constants.Range("C14").Value = "Path for output files. Could be relative or full path."

# Row 15 keeps the wrap-text style used elsewhere in column C, with no value.
$constants.Range("C15").WrapText = $true

# Extend the used range down to row 989 (two new blank formatted rows).
$constants.Rows.Item(988).RowHeight = 14.25
$constants.Rows.Item(989).RowHeight = 14.25

# Move the active selection to C6.
$constants.Activate()
$constants.Range("C6").Select()

# Leave "Settings" as the active/selected sheet, matching the final state.
$settings.Activate()
